$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing data rows (2 onward) contents while preserving row1 header formatting
$ws.Range("A2:T17").ClearContents()

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Il17c"
$ws.Range("C2").Value = "Il17ra"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.784813666666667
$ws.Range("H2").Value = 5.354441
$ws.Range("I2").Value = 0.3928052716376136
$ws.Range("J2").Value = 0.3928052716376136
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 17.37362766666667
$ws.Range("N2").Value = 52.120883
$ws.Range("O2").Value = 0.4119962501387955
$ws.Range("P2").Value = 0.4119962501387954
$ws.Range("Q2").Value = 31.00868809904478
$ws.Range("R2").Value = 279.078192891403
$ws.Range("S2").Value = 0.1618342989494478
$ws.Range("T2").Value = 0.1618342989494477

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Il17c"
$ws.Range("C3").Value = "Il17ra"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 1.784813666666667
$ws.Range("H3").Value = 5.354441
$ws.Range("I3").Value = 0.3928052716376136
$ws.Range("J3").Value = 0.3928052716376136
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 7.318911
$ws.Range("N3").Value = 21.956733
$ws.Range("O3").Value = 0.1735598313117363
$ws.Range("P3").Value = 0.1735598313117363
$ws.Range("Q3").Value = 13.062892377917
$ws.Range("R3").Value = 117.566031401253
$ws.Range("S3").Value = 0.06817521668378497
$ws.Range("T3").Value = 0.06817521668378497

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Il17c"
$ws.Range("C4").Value = "Il17ra"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.784813666666667
$ws.Range("H4").Value = 5.354441
$ws.Range("I4").Value = 0.3928052716376136
$ws.Range("J4").Value = 0.3928052716376136
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 13.84501566666667
$ws.Range("N4").Value = 41.535047
$ws.Range("O4").Value = 0.3283191425083613
$ws.Range("P4").Value = 0.3283191425083613
$ws.Range("Q4").Value = 24.71077317708077
$ws.Range("R4").Value = 222.396958593727
$ws.Range("S4").Value = 0.1289654899568252
$ws.Range("T4").Value = 0.1289654899568252

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Il17c"
$ws.Range("C5").Value = "Il17ra"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 1.784813666666667
$ws.Range("H5").Value = 5.354441
$ws.Range("I5").Value = 0.3928052716376136
$ws.Range("J5").Value = 0.3928052716376136
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.631828666666667
$ws.Range("N5").Value = 10.895486
$ws.Range("O5").Value = 0.08612477604110705
$ws.Range("P5").Value = 0.08612477604110705
$ws.Range("Q5").Value = 6.482137439258445
$ws.Range("R5").Value = 58.33923695332599
$ws.Range("S5").Value = 0.03383026604755569
$ws.Range("T5").Value = 0.03383026604755569

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Il17c"
$ws.Range("C6").Value = "Il17ra"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.5783573333333334
$ws.Range("H6").Value = 1.735072
$ws.Range("I6").Value = 0.1272860095518501
$ws.Range("J6").Value = 0.1272860095518501
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 17.37362766666667
$ws.Range("N6").Value = 52.120883
$ws.Range("O6").Value = 0.4119962501387955
$ws.Range("P6").Value = 0.4119962501387954
$ws.Range("Q6").Value = 10.04816496761956
$ws.Range("R6").Value = 90.433484708576
$ws.Range("S6").Value = 0.05244135863049313
$ws.Range("T6").Value = 0.05244135863049312

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Il17c"
$ws.Range("C7").Value = "Il17ra"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.5783573333333334
$ws.Range("H7").Value = 1.735072
$ws.Range("I7").Value = 0.1272860095518501
$ws.Range("J7").Value = 0.1272860095518501
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 7.318911
$ws.Range("N7").Value = 21.956733
$ws.Range("O7").Value = 0.1735598313117363
$ws.Range("P7").Value = 0.1735598313117363
$ws.Range("Q7").Value = 4.232945848864
$ws.Range("R7").Value = 38.096512639776
$ws.Range("S7").Value = 0.02209173834616315
$ws.Range("T7").Value = 0.02209173834616315

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Il17c"
$ws.Range("C8").Value = "Il17ra"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.5783573333333334
$ws.Range("H8").Value = 1.735072
$ws.Range("I8").Value = 0.1272860095518501
$ws.Range("J8").Value = 0.1272860095518501
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 13.84501566666667
$ws.Range("N8").Value = 41.535047
$ws.Range("O8").Value = 0.3283191425083613
$ws.Range("P8").Value = 0.3283191425083613
$ws.Range("Q8").Value = 8.007366340931556
$ws.Range("R8").Value = 72.066297068384
$ws.Range("S8").Value = 0.0417904335093745
$ws.Range("T8").Value = 0.0417904335093745

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Il17c"
$ws.Range("C9").Value = "Il17ra"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.5783573333333334
$ws.Range("H9").Value = 1.735072
$ws.Range("I9").Value = 0.1272860095518501
$ws.Range("J9").Value = 0.1272860095518501
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 3.631828666666667
$ws.Range("N9").Value = 10.895486
$ws.Range("O9").Value = 0.08612477604110705
$ws.Range("P9").Value = 0.08612477604110705
$ws.Range("Q9").Value = 2.100494742776889
$ws.Range("R9").Value = 18.904452684992
$ws.Range("S9").Value = 0.0109624790658193
$ws.Range("T9").Value = 0.0109624790658193

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Il17c"
$ws.Range("C10").Value = "Il17ra"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.712277333333333
$ws.Range("H10").Value = 5.136832
$ws.Range("I10").Value = 0.3768413339724513
$ws.Range("J10").Value = 0.3768413339724513
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 17.37362766666667
$ws.Range("N10").Value = 52.120883
$ws.Range("O10").Value = 0.4119962501387955
$ws.Range("P10").Value = 0.4119962501387954
$ws.Range("Q10").Value = 29.74846885140623
$ws.Range("R10").Value = 267.736219662656
$ws.Range("S10").Value = 0.1552572164939514
$ws.Range("T10").Value = 0.1552572164939514

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Il17c"
$ws.Range("C11").Value = "Il17ra"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.712277333333333
$ws.Range("H11").Value = 5.136832
$ws.Range("I11").Value = 0.3768413339724513
$ws.Range("J11").Value = 0.3768413339724513
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 7.318911
$ws.Range("N11").Value = 21.956733
$ws.Range("O11").Value = 0.1735598313117363
$ws.Range("P11").Value = 0.1735598313117363
$ws.Range("Q11").Value = 12.532005409984
$ws.Range("R11").Value = 112.788048689856
$ws.Range("S11").Value = 0.06540451835554834
$ws.Range("T11").Value = 0.06540451835554832

# Row 12
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Il17c"
$ws.Range("C12").Value = "Il17ra"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1.712277333333333
$ws.Range("H12").Value = 5.136832
$ws.Range("I12").Value = 0.3768413339724513
$ws.Range("J12").Value = 0.3768413339724513
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 13.84501566666667
$ws.Range("N12").Value = 41.535047
$ws.Range("O12").Value = 0.3283191425083613
$ws.Range("P12").Value = 0.3283191425083613
$ws.Range("Q12").Value = 23.70650650567822
$ws.Range("R12").Value = 213.358558551104
$ws.Range("S12").Value = 0.1237242236315422
$ws.Range("T12").Value = 0.1237242236315422

# Row 13
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Il17c"
$ws.Range("C13").Value = "Il17ra"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1.712277333333333
$ws.Range("H13").Value = 5.136832
$ws.Range("I13").Value = 0.3768413339724513
$ws.Range("J13").Value = 0.3768413339724513
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 3.631828666666667
$ws.Range("N13").Value = 10.895486
$ws.Range("O13").Value = 0.08612477604110705
$ws.Range("P13").Value = 0.08612477604110705
$ws.Range("Q13").Value = 6.218697904483556
$ws.Range("R13").Value = 55.968281140352
$ws.Range("S13").Value = 0.03245537549140939
$ws.Range("T13").Value = 0.03245537549140939

# Row 14
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Il17c"
$ws.Range("C14").Value = "Il17ra"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.4683136666666667
$ws.Range("H14").Value = 1.404941
$ws.Range("I14").Value = 0.103067384838085
$ws.Range("J14").Value = 0.103067384838085
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 17.37362766666667
$ws.Range("N14").Value = 52.120883
$ws.Range("O14").Value = 0.4119962501387955
$ws.Range("P14").Value = 0.4119962501387954
$ws.Range("Q14").Value = 8.136307275878112
$ws.Range("R14").Value = 73.226765482903
$ws.Range("S14").Value = 0.04246337606490316
$ws.Range("T14").Value = 0.04246337606490314

# Row 15
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Il17c"
$ws.Range("C15").Value = "Il17ra"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.4683136666666667
$ws.Range("H15").Value = 1.404941
$ws.Range("I15").Value = 0.103067384838085
$ws.Range("J15").Value = 0.103067384838085
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 7.318911
$ws.Range("N15").Value = 21.956733
$ws.Range("O15").Value = 0.1735598313117363
$ws.Range("P15").Value = 0.1735598313117363
$ws.Range("Q15").Value = 3.427546046417
$ws.Range("R15").Value = 30.847914417753
$ws.Range("S15").Value = 0.01788835792623984
$ws.Range("T15").Value = 0.01788835792623983

# Row 16
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Il17c"
$ws.Range("C16").Value = "Il17ra"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.4683136666666667
$ws.Range("H16").Value = 1.404941
$ws.Range("I16").Value = 0.103067384838085
$ws.Range("J16").Value = 0.103067384838085
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 13.84501566666667
$ws.Range("N16").Value = 41.535047
$ws.Range("O16").Value = 0.3283191425083613
$ws.Range("P16").Value = 0.3283191425083613
$ws.Range("Q16").Value = 6.483810051914111
$ws.Range("R16").Value = 58.354290467227
$ws.Range("S16").Value = 0.03383899541061934
$ws.Range("T16").Value = 0.03383899541061933

# Row 17
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Il17c"
$ws.Range("C17").Value = "Il17ra"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.4683136666666667
$ws.Range("H17").Value = 1.404941
$ws.Range("I17").Value = 0.103067384838085
$ws.Range("J17").Value = 0.103067384838085
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 3.631828666666667
$ws.Range("N17").Value = 10.895486
$ws.Range("O17").Value = 0.08612477604110705
$ws.Range("P17").Value = 0.08612477604110705
$ws.Range("Q17").Value = 1.700834999591778
$ws.Range("R17").Value = 15.307514996326
$ws.Range("S17").Value = 0.008876655436322661
$ws.Range("T17").Value = 0.008876655436322661

Write-Output "done"